# Updates the cryptos list (Price column D, Volume(1h) column E) with the
# latest scraped values, per the "Updated cryptos list ... with GitHub
# Actions" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "27.096.05"; E = "  -2.37%  " },
    @{ Row = 3; D = "1.873.54"; E = "  -1.64%  " },
    @{ Row = 4; D = "1.001"; E = "  +0.23%  " },
    @{ Row = 5; D = "306.64"; E = "  -1.60%  " },
    @{ Row = 6; D = "1.001"; E = "  +0.13%  " },
    @{ Row = 7; D = "0.5089"; E = "  -1.77%  " },
    @{ Row = 8; D = "0.3727"; E = "  -1.86%  " },
    @{ Row = 9; D = "0.07176"; E = "  -1.05%  " },
    @{ Row = 10; D = "0.8929"; E = "  -1.24%  " },
    @{ Row = 11; D = "20.82"; E = "  -2.33%  " },
    @{ Row = 12; D = "1.892.82"; E = "  -0.89%  " },
    @{ Row = 13; D = "0.07519"; E = "  -1.83%  " },
    @{ Row = 14; D = "5.300"; E = "  -2.83%  " },
    @{ Row = 15; D = "90.81"; E = "  -1.86%  " },
    @{ Row = 16; D = "1.001"; E = "  +0.23%  " },
    @{ Row = 17; D = "0.000008669"; E = "  -0.46%  " },
    @{ Row = 18; D = "14.08"; E = "  -2.65%  " },
    @{ Row = 19; D = "1.001"; E = "  +0.13%  " },
    @{ Row = 20; D = "27.139.15"; E = "  -2.37%  " },
    @{ Row = 21; D = "5.023"; E = "  -2.47%  " },
    @{ Row = 22; D = "2.123.56"; E = "  -3.25%  " },
    @{ Row = 23; D = "10.41"; E = "  -4.11%  " },
    @{ Row = 24; D = "6.483"; E = "  -2.37%  " },
    @{ Row = 25; D = "1.831"; E = "  -1.37%  " },
    @{ Row = 26; D = "146.52"; E = "  -4.40%  " },
    @{ Row = 27; D = "18.00"; E = "  -1.76%  " },
    @{ Row = 28; D = "2.073"; E = "  -4.52%  " },
    @{ Row = 29; D = "113.14"; E = "  -0.83%  " },
    @{ Row = 30; D = "4.661"; E = "  -3.32%  " },
    @{ Row = 31; D = "4.694"; E = "  -3.18%  " },
    @{ Row = 32; D = "0.09231"; E = "  +1.50%  " },
    @{ Row = 33; D = "0.05118"; E = "  -3.31%  " },
    @{ Row = 34; D = "3.094"; E = "  -2.74%  " },
    @{ Row = 35; D = "1.155"; E = "  -6.26%  " },
    @{ Row = 36; D = "0.7262"; E = "  -7.43%  " },
    @{ Row = 37; D = "3.173"; E = "  +3.04%  " },
    @{ Row = 38; D = "0.02028"; E = "  -3.06%  " },
    @{ Row = 39; D = "2.496"; E = "  -4.18%  " },
    @{ Row = 40; D = $null; E = "  -1.60%  " },
    @{ Row = 41; D = "0.5306"; E = "  -5.21%  " },
    @{ Row = 42; D = "6.513"; E = "  -2.95%  " },
    @{ Row = 43; D = "116.83"; E = "  -0.85%  " },
    @{ Row = 44; D = "8.361"; E = "  -2.63%  " },
    @{ Row = 45; D = $null; E = "  -2.90%  " },
    @{ Row = 46; D = "0.4641"; E = "  -4.23%  " },
    @{ Row = 47; D = "1.001"; E = "  +0.15%  " },
    @{ Row = 48; D = "10.00"; E = "  -4.75%  " },
    @{ Row = 49; D = "1.564"; E = "  -3.31%  " },
    @{ Row = 50; D = "36.96"; E = "  -0.44%  " },
    @{ Row = 51; D = "63.69"; E = "  -4.85%  " }
)

foreach ($item in $updates) {
    if ($null -ne $item.D) {
        $dCell = $ws.Range("D" + $item.Row)
        # Force the cell to remain text so numeric-looking price strings
        # (e.g. "1.001", "306.64") aren't silently coerced into numbers
        # (which would also strip significant trailing zeros, e.g. "5.300").
        $dCell.NumberFormat = "@"
        $dCell.Value = $item.D
        $dCell.NumberFormat = "General"
        $dCell.Style = "Normal"
    }

    $eCell = $ws.Range("E" + $item.Row)
    $eCell.Value = $item.E
}
